$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows 8 & 9: shift fiscal-period / publish-date labels left, add new final column ---
$ws.Range("D8").Value = $ws.Range("E8").Value2
$ws.Range("E8").Value = $ws.Range("F8").Value2
$ws.Range("F8").Value = $ws.Range("G8").Value2
$ws.Range("G8").Value = $ws.Range("H8").Value2
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

$ws.Range("D9").Value = $ws.Range("E9").Value2
$ws.Range("E9").Value = $ws.Range("F9").Value2
$ws.Range("F9").Value = $ws.Range("G9").Value2
$ws.Range("G9").Value = "1402-02-30 (7)"
$ws.Range("H9").Value = "1402-02-30"

# --- Data rows 12-58: shift each period's figures one column left (drop oldest year),
#     then place the newly reported fiscal-year figures into column H ---
$r = 12
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 1825948

$r = 13
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 1330

$r = 14
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 3035359

$r = 15
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 4610307

$r = 16
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 2029647

$r = 17
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 18
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 11502591

$r = 19
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 20274

$r = 20
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 21
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 22
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 1425095

$r = 23
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 193

$r = 24
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 25
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 26
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 1445882

$r = 27
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 12948473

$r = 29
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 1922587

$r = 30
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 31
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 3248690

$r = 32
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 412336

$r = 33
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 9946

$r = 34
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 0

$r = 35
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 36
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 37
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 5593559

$r = 38
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 135044

$r = 39
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 40
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 41
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 853800

$r = 42
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 988844

$r = 43
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 6582403

$r = 45
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 1826000

$r = 46
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 47
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 0

$r = 48
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 49
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 0

$r = 50
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 182600

$r = 51
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 0

$r = 52
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 53
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 54
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 55
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2

$r = 56
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 4357470

$r = 57
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 6366070

$r = 58
$ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 5).Value2
$ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 6).Value2
$ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 7).Value2
$ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 8).Value2
$ws.Cells.Item($r, 8).Value = 12948473
